# Atualização de bases das ligas, do dia: 25-05-2024 às 15:10
#
# The underlying data rows (columns B:AB) for several fixtures got
# re-shuffled between rows that share the same match date (column D /
# column A index stay put - only the match data itself moved rows).
# Snapshot every affected row's B:AB values first, then write them back
# out to their new homes so rows that are part of a multi-row rotation
# (not just a simple two-row swap) are handled correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All rows whose B:AB contents are being relocated.
$rows = @(18,19,42,43,81,82,90,91,94,95,98,99,100,101,102,103,108,109,118,119,122,123,124,126,127,131,132,133)

# Snapshot current contents (B:AB) of every affected row before any writes.
$data = @{}
foreach ($r in $rows) {
    $range = "B" + $r + ":AB" + $r
    $data[$r] = $ws.Range($range).Value2
}

# perm[destinationRow] = sourceRow whose snapshotted content now belongs
# at destinationRow. Pairs are simple swaps; rows 122/123/127 and
# 131/132/133 form 3-way rotations.
$perm = @{
    18=19;   19=18;
    42=43;   43=42;
    81=82;   82=81;
    90=91;   91=90;
    94=95;   95=94;
    98=99;   99=98;
    100=101; 101=100;
    102=103; 103=102;
    108=109; 109=108;
    118=119; 119=118;
    122=123; 123=127; 127=122;
    124=126; 126=124;
    131=133; 133=132; 132=131;
}

foreach ($r in $rows) {
    $src = $perm[$r]
    $range = "B" + $r + ":AB" + $r
    $ws.Range($range).Value2 = $data[$src]
}
